$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 44) into the
# new row 45, then fill in the new block's values.
$ws.Range("A44:E44").Copy()
$ws.Range("A45:E45").PasteSpecial(-4122)

$ws.Range("A45").Value = "charcoal_pit:log_pile"
$ws.Range("B45").Value = 300
$ws.Range("C45").Value = 3

# Booleans: copying format straight from an already-boolean cell (D44/E44)
# collapses onto a different (but visually identical) style record. Seed the
# new boolean cells' formatting from a same-block non-boolean cell instead
# (B44, style shared with C44) so the xf index matches, then set the values.
$ws.Range("B44").Copy()
$ws.Range("D45:E45").PasteSpecial(-4122)

$ws.Range("D45").Value = $false
$ws.Range("E45").Value = $true
